$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(13, 0.5,  3081421.847401,  38928.218222,  3042493.629179,  25787.17308033333, 1140010.440327, 20161.456375, 1119848.983952, 27700.59754266667),
    @(14, 0.55, 3403300.513984,  42702.714682,  3360597.799302,  29979.56488633333, 1176120.188236, 20759.429981, 1155360.758255, 32060.995423),
    @(15, 0.6,  4156758.015218,  44285.039584,  4112472.975634,  33080.04253033333, 1358508.266211, 20728.391044, 1337779.875167, 35302.04909766667),
    @(16, 0.65, 4656527.744997,  44774.934837,  4611752.81016,   35804.18152733333, 1447754.652049, 20574.237238, 1427180.414811, 38187.99589733333)
)

foreach ($row in $data) {
    $r = $row[0]
    for ($col = 1; $col -le 9; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}
